# Word COM-interop edit script
# Applies the "dark energy" -> "chemistry" content swap described by the
# commit diff, including the font-name fix (TimesNewToman -> Times New
# Roman) everywhere in the document.

$d = $word.ActiveDocument
$BR = [char]11

function Replace-WholeWord($findText, $replaceText) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute($findText, $true, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null
}

function Set-ParagraphSegments($para, [object[]]$segments) {
    # Concatenate the target segments into one string and push it into the
    # paragraph's range (excluding the trailing paragraph mark). Word merges
    # same-formatted runs automatically, so we restore the desired run
    # boundaries afterwards by toggling the font at each boundary, which
    # forces the engine to re-split the run without altering the visible
    # formatting.
    $full = [string]::Join('', $segments)
    $pStart = $para.Range.Start
    $pEnd = $para.Range.End - 1
    $target = $d.Range($pStart, $pEnd)
    $origFont = $target.Font.Name
    $target.Text = $full

    $offset = $pStart
    $boundaries = @()
    for ($i = 0; $i -lt $segments.Length - 1; $i++) {
        $offset += $segments[$i].Length
        $boundaries += $offset
    }

    $paraEndNow = $para.Range.End - 1
    foreach ($b in $boundaries) {
        $tail = $d.Range($b, $paraEndNow)
        $tail.Font.Name = "Arial"
        $tail.Font.Name = $origFont
    }
}

# ---------------------------------------------------------------------
# 1. Font-name fix across the whole document: TimesNewToman -> Times New
#    Roman (a formatting-only change, so it does not disturb run layout).
# ---------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Font.Name = "TimesNewToman"
$find.Replacement.ClearFormatting()
$find.Replacement.Font.Name = "Times New Roman"
$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $true, $null, 2) | Out-Null

# ---------------------------------------------------------------------
# 2. Title / author / e-mail address.
# ---------------------------------------------------------------------
Replace-WholeWord "Disentangling the Enigma of Dark Energy" "Delving into the Realm of Chemistry: Unveiling the Wonders at the Molecular Level"
Replace-WholeWord "Penelope Williams" "Alice Patterson"
Replace-WholeWord "pwilliams@amail" "apatterson@eduvate"
Replace-WholeWord "com" "org"

# ---------------------------------------------------------------------
# 3. Main body paragraph (dark-energy essay -> chemistry essay).
# ---------------------------------------------------------------------
$para4Segs = @(
  'Chemistry serves as an intriguing discipline that delves into the fundamental principles and interactions underlying the properties, composition, and transformations of matter',
  '.',
  ' It acts as a magic key that opens the gates to the secrets hidden within substances, and it holds the power to unravel the mysteries that dictate how they behave',
  '.',
  ' With fascination, we delve into the microscopic world, where atoms, the building blocks of all substances, engage in captivating dances of interactions',
  '.',
  ' In our journey of discovery, we perceive the mesmerizing interplay of particles as they exchange energy and form new combinations, thus weaving the intricate tapestry of the material world',
  '.',
  ($BR),
  ($BR + 'Chemistry not only enables us to decipher the mysteries of matter but also empowers us to harness its potential for the betterment of society'),
  '.',
  ' Through its lens, we uncover innovative solutions to global challenges, such as the development of cleaner energy sources, the engineering of advanced materials, and the synthesis of life-saving pharmaceuticals',
  '.',
  ' Chemistry serves as an indispensable tool, aiding us in safeguarding the environment, enhancing human health, and facilitating technological advancements that shape our modern world',
  '.',
  ' As we unveil the intricacies of chemistry, we unlock the potential to forge a sustainable and prosperous future',
  '.',
  ($BR),
  ($BR + 'The versatility of chemistry extends to its diverse applications across various fields'),
  '.',
  ' It plays a profound role in the medical realm, facilitating the development of effective medications and therapies',
  '.',
  ' In agriculture, chemistry contributes to developing more productive crop varieties and devising innovative pest management strategies',
  '.',
  ' Furthermore, it finds its place in materials science, leading to the creation of advanced materials with tailored properties that serve a multitude of purposes',
  '.',
  ' Chemistry''s impact is apparent in energy production, propelling the transition to sustainable and efficient energy sources',
  '.',
  ' The footprints of chemistry are ubiquitous in our daily lives beyond these core areas; from the clothes we wear, to the food we consume, to the devices we utilize, its presence is pervasive, making it a field of endless exploration and discovery',
  '.'
)

$bodyPara = $d.Paragraphs.Item(5)
Set-ParagraphSegments $bodyPara $para4Segs

# ---------------------------------------------------------------------
# 4. Summary paragraph.
# ---------------------------------------------------------------------
$para6Segs = @(
  'Chemistry stands as a testament to the wonders of the microscopic world, offering boundless opportunities to explore the intricacies of matter and its transformations',
  '.',
  ' Its influence transcends the boundaries of scientific inquiry; it impacts various areas of human endeavor, including medicine, energy, agriculture, and materials science',
  '.',
  ' By delving into the realm of chemistry, we cultivate critical thinking abilities, nurture our curiosity, and gain a deeper understanding of the world around us',
  '.',
  ' It is this pursuit of knowledge, coupled with an inquisitive spirit, that enables us to unravel the secrets of matter and harness its potential for societal progress',
  '.'
)

$summaryPara = $d.Paragraphs.Item(7)
Set-ParagraphSegments $summaryPara $para6Segs

# Insert the <w:lastRenderedPageBreak/> marker at the very start of the
# Summary paragraph's first run.
$summaryPara.Range.Characters.Item(1).InsertXML('<w:lastRenderedPageBreak xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>')

# ---------------------------------------------------------------------
# 5. A new, empty trailing paragraph at the very end of the document body.
# ---------------------------------------------------------------------
$d.Paragraphs.Item($d.Paragraphs.Count).Range.InsertParagraphAfter() | Out-Null
